# Apply "added towns and some more latitude/longitude" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Office Location (column L) additions: town names ---
$ws.Range("L2").Value  = "Alexandria"
$ws.Range("L3").Value  = "Richmond"
$ws.Range("L4").Value  = "Richmond"
$ws.Range("L5").Value  = "Ashburn"
$ws.Range("L6").Value  = "Ashburn"
$ws.Range("L7").Value  = "Ashburn"
$ws.Range("L9").Value  = "Online"
$ws.Range("L10").Value = "Leesburg"
$ws.Range("L11").Value = "Leesburg"
$ws.Range("L12").Value = "Vienna"
$ws.Range("L24").Value = "Pittsburgh"
$ws.Range("L25").Value = "Pittsburgh"
$ws.Range("L26").Value = "McKeesport"
$ws.Range("L27").Value = "Pittsburgh"
$ws.Range("L28").Value = "Pittsburgh and McKeesport"
$ws.Range("L29").Value = "Pittsburgh"
$ws.Range("L30").Value = "Pittsburgh"
$ws.Range("L31").Value = "Pittsburgh"

# --- Latitude / Longitude (columns J and K) additions ---
$ws.Range("J17").Value = 39.114571664302197
$ws.Range("K17").Value = -77.5405475733521

$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0

# --- Row height adjustments (rows re-wrap after the new numeric lat/long values) ---
$ws.Rows.Item(20).RowHeight = 84
$ws.Rows.Item(21).RowHeight = 63

# --- View-state refresh to match the new selection/scroll location ---
$ws.Range("M28").Select()
